$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: measure/dimension labels - rename dimension:aragon+sexo into a single
# sdmx-dimension:refArea / measure:sexo pairing, and edad-grupos-quinquenales
# becomes a "measure" instead of "dimension"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "iaest-measure:edad-grupos-quinquenales"

# Row 3: medida/dim markers follow the same re-shuffle
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "medida"

# Row 4: types
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-Provincia"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "xsd:int"

# Row 5 no longer exists - clear it entirely
$ws.Range("A5:G5").Clear()
